$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2791
$ws.Range("J17").Value = 2418.6128
$ws.Range("L17").Value = 7255.8384
$ws.Range("N17").Value = -7591.8384
$ws.Range("H19").Value = 527930.25
$ws.Range("I19").Value = 1250598
$ws.Range("K19").Value = 1250598
$ws.Range("M19").Value = -1250423
$ws.Range("H28").Value = 3342.5715
$ws.Range("I28").Value = 233.16667
$ws.Range("J28").Value = 21999
$ws.Range("K28").Value = 233.16667
$ws.Range("L28").Value = 21999
$ws.Range("M28").Value = 251.83333
$ws.Range("N28").Value = -22969
$ws.Range("H32").Value = 2685
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2685
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2685
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3337
$ws.Range("H43").Value = 1416.0769
$ws.Range("I43").Value = 1298.6666
$ws.Range("K43").Value = 1298.6666
$ws.Range("M43").Value = -1229.6666
$ws.Range("H69").Value = 2999.5
$ws.Range("I69").Value = 2999.5
$ws.Range("K69").Value = 8998.5
$ws.Range("M69").Value = -8124.5
$ws.Range("H72").Value = 2999.5
$ws.Range("I72").Value = 2999.5
$ws.Range("K72").Value = 26995.5
$ws.Range("M72").Value = -22627.5
$ws.Range("H88").Value = 2485
$ws.Range("I88").Value = 1500
$ws.Range("K88").Value = 1500
$ws.Range("M88").Value = -1094
$ws.Range("H91").Value = 2485
$ws.Range("I91").Value = 1500
$ws.Range("K91").Value = 1500
$ws.Range("M91").Value = -96
$ws.Range("H97").Value = 938.4286
$ws.Range("J97").Value = 946.8333
$ws.Range("L97").Value = 2840.4999
$ws.Range("N97").Value = -3832.4999
$ws.Range("H98").Value = 2841.9333
$ws.Range("I98").Value = 2949.4546
$ws.Range("J98").Value = 2546.25
$ws.Range("K98").Value = 2949.4546
$ws.Range("L98").Value = 2546.25
$ws.Range("M98").Value = -1451.4546
$ws.Range("N98").Value = -5542.25
$ws.Range("H122").Value = 2841.9333
$ws.Range("I122").Value = 2949.4546
$ws.Range("J122").Value = 2546.25
$ws.Range("K122").Value = 8848.363799999999
$ws.Range("L122").Value = 7638.75
$ws.Range("M122").Value = -6398.363799999999
$ws.Range("N122").Value = -12538.75
$ws.Range("H132").Value = 1329
$ws.Range("J132").Value = 1300.8572
$ws.Range("L132").Value = 3902.5716
$ws.Range("N132").Value = -8962.571599999999
$ws.Range("H137").Value = 1400.75
$ws.Range("I137").Value = 1100.6666
$ws.Range("K137").Value = 3301.9998
$ws.Range("M137").Value = -751.9998000000001
$ws.Range("H138").Value = 3171.0889
$ws.Range("I138").Value = 4460.25
$ws.Range("J138").Value = 2459.8276
$ws.Range("K138").Value = 13380.75
$ws.Range("L138").Value = 7379.4828
$ws.Range("M138").Value = -8240.75
$ws.Range("N138").Value = -17659.4828
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 25000250
$ws.Range("I13").Value = 25000250
$ws.Range("K13").Value = 25000250
$ws.Range("M13").Value = -25000106
$ws.Range("H61").Value = 6203.591
$ws.Range("I61").Value = 6328.2354
$ws.Range("K61").Value = 6328.2354
$ws.Range("M61").Value = -6116.2354
$ws.Range("H74").Value = 1611.75
$ws.Range("I74").Value = 482.33334
$ws.Range("K74").Value = 482.33334
$ws.Range("M74").Value = 391.66666
$ws.Range("H77").Value = 1611.75
$ws.Range("I77").Value = 482.33334
$ws.Range("K77").Value = 2411.6667
$ws.Range("M77").Value = 1956.3333
$ws.Range("H136").Value = 6203.591
$ws.Range("I136").Value = 6328.2354
$ws.Range("K136").Value = 18984.7062
$ws.Range("M136").Value = -16434.7062
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 48000
$ws.Range("J68").Value = 48000
$ws.Range("L68").Value = 48000
$ws.Range("N68").Value = -49622
$ws.Range("H71").Value = 48000
$ws.Range("J71").Value = 48000
$ws.Range("L71").Value = 144000
$ws.Range("N71").Value = -152112
$ws.Range("H107").Value = 940.1429000000001
$ws.Range("I107").Value = 876.4
$ws.Range("J107").Value = 1099.5
$ws.Range("K107").Value = 876.4
$ws.Range("L107").Value = 1099.5
$ws.Range("M107").Value = 1043.6
$ws.Range("N107").Value = -4939.5
$ws.Range("H134").Value = 4410.9697
$ws.Range("I134").Value = 4614.207
$ws.Range("K134").Value = 13842.621
$ws.Range("M134").Value = -11307.621
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2614.2
$ws.Range("I31").Value = 2240
$ws.Range("K31").Value = 2240
$ws.Range("M31").Value = -1945
$ws.Range("H34").Value = 2614.2
$ws.Range("I34").Value = 2240
$ws.Range("K34").Value = 2240
$ws.Range("M34").Value = -2038
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184
$ws.Range("H105").Value = 831.61536
$ws.Range("I105").Value = 808.5
$ws.Range("K105").Value = 808.5
$ws.Range("M105").Value = 938.5
$ws.Range("H107").Value = 1512.875
$ws.Range("I107").Value = 1131.8
$ws.Range("J107").Value = 2148
$ws.Range("K107").Value = 1131.8
$ws.Range("L107").Value = 2148
$ws.Range("M107").Value = 788.2
$ws.Range("N107").Value = -5988
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 644.1429000000001
$ws.Range("I11").Value = 585
$ws.Range("K11").Value = 1755
$ws.Range("M11").Value = -1615
$ws.Range("H140").Value = 1849.1875
$ws.Range("I140").Value = 1029.9231
$ws.Range("K140").Value = 3089.7693
$ws.Range("M140").Value = 2090.2307
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 130.83333
$ws.Range("I2").Value = 62
$ws.Range("J2").Value = 165.25
$ws.Range("K2").Value = 62
$ws.Range("L2").Value = 165.25
$ws.Range("M2").Value = 51
$ws.Range("N2").Value = -391.25
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H102").Value = 2550.75
$ws.Range("I102").Value = 2145.1428
$ws.Range("J102").Value = 3497.1667
$ws.Range("K102").Value = 2145.1428
$ws.Range("L102").Value = 3497.1667
$ws.Range("M102").Value = -523.1428000000001
$ws.Range("N102").Value = -6741.1667
$ws.Range("H107").Value = 147.5
$ws.Range("I107").Value = 130
$ws.Range("K107").Value = 130
$ws.Range("M107").Value = 1790
$ws.Range("H132").Value = 1674914.4
$ws.Range("I132").Value = 2138624.2
$ws.Range("J132").Value = 5559.4
$ws.Range("K132").Value = 6415872.600000001
$ws.Range("L132").Value = 16678.2
$ws.Range("M132").Value = -6413342.600000001
$ws.Range("N132").Value = -21738.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2476.6667
$ws.Range("I7").Value = 2762.5557
$ws.Range("J7").Value = 2047.8334
$ws.Range("K7").Value = 2762.5557
$ws.Range("L7").Value = 2047.8334
$ws.Range("M7").Value = -2650.5557
$ws.Range("N7").Value = -2271.8334
$ws.Range("H17").Value = 499.85715
$ws.Range("I17").Value = 499.85715
$ws.Range("K17").Value = 499.85715
$ws.Range("M17").Value = -329.85715
$ws.Range("H126").Value = 2476.6667
$ws.Range("I126").Value = 2762.5557
$ws.Range("J126").Value = 2047.8334
$ws.Range("K126").Value = 8287.667099999999
$ws.Range("L126").Value = 6143.5002
$ws.Range("M126").Value = -5817.667099999999
$ws.Range("N126").Value = -11083.5002
$ws.Range("H136").Value = 2010.5
$ws.Range("I136").Value = 1644.9231
$ws.Range("J136").Value = 2961
$ws.Range("K136").Value = 4934.7693
$ws.Range("L136").Value = 8883
$ws.Range("M136").Value = -2384.7693
$ws.Range("N136").Value = -13983
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 16000
$ws.Range("I20").Value = 12000
$ws.Range("J20").Value = 18000
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = -11760
$ws.Range("N20").Value = -18480
$ws.Range("H82").Value = 45000
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 45000
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 45000
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 45000
$ws.Range("N85").Value = -47652
$ws.Range("H107").Value = 644.45
$ws.Range("J107").Value = 876.8333
$ws.Range("L107").Value = 2630.4999
$ws.Range("N107").Value = -6470.4999
$ws.Range("H132").Value = 1841.75
$ws.Range("I132").Value = 1365.1818
$ws.Range("J132").Value = 3589.1667
$ws.Range("K132").Value = 4095.5454
$ws.Range("L132").Value = 10767.5001
$ws.Range("M132").Value = -1565.5454
$ws.Range("N132").Value = -15827.5001
$ws.Range("H136").Value = 32681490
$ws.Range("I136").Value = 50506416
$ws.Range("J136").Value = 2460
$ws.Range("K136").Value = 151519248
$ws.Range("L136").Value = 7380
$ws.Range("M136").Value = -151516698
$ws.Range("N136").Value = -12480
